{"js": "// Replace the three-digit \u00f7 one-digit division problems in the table\n// with the updated values from the commit. Each old value is unique in\n// the document, so an exact-match search safely targets the right cell.\nconst replacements = [\n  [\"396\u00f74=99, 0\", \"869\u00f76=144, 5\"],\n  [\"274\u00f72=137, 0\", \"236\u00f76=39, 2\"],\n  [\"360\u00f74=90, 0\", \"355\u00f78=44, 3\"],\n  [\"921\u00f79=102, 3\", \"994\u00f72=497, 0\"],\n  [\"315\u00f77=45, 0\", \"533\u00f75=106, 3\"],\n  [\"734\u00f76=122, 2\", \"328\u00f77=46, 6\"],\n  [\"582\u00f78=72, 6\", \"195\u00f75=39, 0\"],\n  [\"515\u00f72=257, 1\", \"948\u00f79=105, 3\"],\n  [\"914\u00f73=304, 2\", \"920\u00f74=230, 0\"],\n  [\"270\u00f75=54, 0\", \"458\u00f72=229, 0\"],\n  [\"557\u00f75=111, 2\", \"452\u00f72=226, 0\"],\n  [\"755\u00f79=83, 8\", \"900\u00f72=450, 0\"],\n  [\"468\u00f72=234, 0\", \"429\u00f78=53, 5\"],\n  [\"220\u00f76=36, 4\", \"107\u00f73=35, 2\"],\n  [\"573\u00f77=81, 6\", \"216\u00f77=30, 6\"],\n  [\"993\u00f72=496, 1\", \"190\u00f75=38, 0\"],\n  [\"930\u00f76=155, 0\", \"410\u00f76=68, 2\"],\n  [\"680\u00f76=113, 2\", \"721\u00f77=103, 0\"],\n  [\"930\u00f74=232, 2\", \"762\u00f79=84, 6\"],\n  [\"478\u00f72=239, 0\", \"238\u00f75=47, 3\"],\n  [\"790\u00f75=158, 0\", \"679\u00f78=84, 7\"],\n  [\"491\u00f78=61, 3\", \"287\u00f75=57, 2\"],\n  [\"723\u00f74=180, 3\", \"752\u00f79=83, 5\"],\n  [\"380\u00f76=63, 2\", \"239\u00f78=29, 7\"],\n  [\"759\u00f78=94, 7\", \"443\u00f78=55, 3\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit \u00f7 one-digit division problems in the table\n# with the updated values from the commit. Each old value is unique in\n# the document, so Find/Execute safely targets the right cell each time.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  ,@(\"396\u00f74=99, 0\", \"869\u00f76=144, 5\")\n  ,@(\"274\u00f72=137, 0\", \"236\u00f76=39, 2\")\n  ,@(\"360\u00f74=90, 0\", \"355\u00f78=44, 3\")\n  ,@(\"921\u00f79=102, 3\", \"994\u00f72=497, 0\")\n  ,@(\"315\u00f77=45, 0\", \"533\u00f75=106, 3\")\n  ,@(\"734\u00f76=122, 2\", \"328\u00f77=46, 6\")\n  ,@(\"582\u00f78=72, 6\", \"195\u00f75=39, 0\")\n  ,@(\"515\u00f72=257, 1\", \"948\u00f79=105, 3\")\n  ,@(\"914\u00f73=304, 2\", \"920\u00f74=230, 0\")\n  ,@(\"270\u00f75=54, 0\", \"458\u00f72=229, 0\")\n  ,@(\"557\u00f75=111, 2\", \"452\u00f72=226, 0\")\n  ,@(\"755\u00f79=83, 8\", \"900\u00f72=450, 0\")\n  ,@(\"468\u00f72=234, 0\", \"429\u00f78=53, 5\")\n  ,@(\"220\u00f76=36, 4\", \"107\u00f73=35, 2\")\n  ,@(\"573\u00f77=81, 6\", \"216\u00f77=30, 6\")\n  ,@(\"993\u00f72=496, 1\", \"190\u00f75=38, 0\")\n  ,@(\"930\u00f76=155, 0\", \"410\u00f76=68, 2\")\n  ,@(\"680\u00f76=113, 2\", \"721\u00f77=103, 0\")\n  ,@(\"930\u00f74=232, 2\", \"762\u00f79=84, 6\")\n  ,@(\"478\u00f72=239, 0\", \"238\u00f75=47, 3\")\n  ,@(\"790\u00f75=158, 0\", \"679\u00f78=84, 7\")\n  ,@(\"491\u00f78=61, 3\", \"287\u00f75=57, 2\")\n  ,@(\"723\u00f74=180, 3\", \"752\u00f79=83, 5\")\n  ,@(\"380\u00f76=63, 2\", \"239\u00f78=29, 7\")\n  ,@(\"759\u00f78=94, 7\", \"443\u00f78=55, 3\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
